$d = $word.ActiveDocument

# The document currently ends with the paragraph describing the "View
# Schedule" use case. We append two new paragraphs after it: a blank
# spacer paragraph, then a paragraph with the new plan text.

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# Move into the newly-created (blank) paragraph and insert another
# paragraph break after it, giving us a blank paragraph followed by an
# empty paragraph ready to hold the new text.
$blankParaIndex = $d.Paragraphs.Count
$blankPara = $d.Paragraphs.Item($blankParaIndex)
$r2 = $blankPara.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# Fill in the final (new) paragraph with the plan text.
$newParaIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newParaIndex)
$r3 = $newPara.Range
$r3.Collapse(0)
$r3.InsertAfter("In order to fully implement this use case, we will first create a complete use case description for it, then implement the main flow, then all the alternative flows that we come up with. We will then do some testing to ensure that the system always behaves as expected for various user actions.")
